$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.267.67'
$ws.Range('E2').Value = '  -2.51%  '
$ws.Range('D3').Value = '1.562.99'
$ws.Range('E3').Value = '  -3.71%  '
$ws.Range('E4').Value = '  -0.28%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '206.86'
$cell.Style = "Normal"
$ws.Range('E6').Value = '  -0.26%  '
$cell = $ws.Range('D7')
$cell.NumberFormat = "@"
$cell.Value = '0.478'
$cell.Style = "Normal"
$ws.Range('E7').Value = '  -4.82%  '
$ws.Range('E8').Value = '  -1.53%  '
$ws.Range('E9').Value = '  -3.25%  '
$cell = $ws.Range('D10')
$cell.NumberFormat = "@"
$cell.Value = '17.76'
$cell.Style = "Normal"
$ws.Range('E10').Value = '  -2.71%  '
$ws.Range('E11').Value = '  -0.75%  '
$ws.Range('D12').Value = '1.779.65'
$ws.Range('E12').Value = '  -3.74%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$cell = $ws.Range('D13')
$cell.NumberFormat = "@"
$cell.Value = '4.00'
$cell.Style = "Normal"
$ws.Range('E13').Value = '  -4.30%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.543.94'
$ws.Range('E14').Value = '  -4.95%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = "@"
$cell.Value = '0.503'
$cell.Style = "Normal"
$ws.Range('E15').Value = '  -3.81%  '
$ws.Range('D16').Value = '25.267.65'
$ws.Range('E16').Value = '  -2.49%  '
$cell = $ws.Range('D17')
$cell.NumberFormat = "@"
$cell.Value = '59.20'
$cell.Style = "Normal"
$ws.Range('E17').Value = '  -3.13%  '
$ws.Range('D18').Value = '0.0₃0710'
$ws.Range('E18').Value = '  -3.11%  '
$ws.Range('E19').Value = '  -0.32%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = "@"
$cell.Value = '185.37'
$cell.Style = "Normal"
$ws.Range('E20').Value = '  -3.34%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = "@"
$cell.Value = '4.11'
$cell.Style = "Normal"
$ws.Range('E21').Value = '  -2.80%  '
$ws.Range('E22').Value = '  -3.22%  '
$cell = $ws.Range('D23')
$cell.NumberFormat = "@"
$cell.Value = '5.87'
$cell.Style = "Normal"
$ws.Range('E23').Value = '  -3.19%  '
$ws.Range('E24').Value = '  -0.26%  '
$ws.Range('E25').Value = '  -3.79%  '
$cell = $ws.Range('D26')
$cell.NumberFormat = "@"
$cell.Value = '139.61'
$cell.Style = "Normal"
$ws.Range('E27').Value = '  -6.58%  '
$cell = $ws.Range('D28')
$cell.NumberFormat = "@"
$cell.Value = '6.45'
$cell.Style = "Normal"
$ws.Range('E28').Value = '  -3.98%  '
$cell = $ws.Range('D29')
$cell.NumberFormat = "@"
$cell.Value = '14.80'
$cell.Style = "Normal"
$ws.Range('E29').Value = '  -2.22%  '
$ws.Range('E30').Value = '  -6.40%  '
$cell = $ws.Range('D31')
$cell.NumberFormat = "@"
$cell.Value = '0.0464'
$cell.Style = "Normal"
$ws.Range('E31').Value = '  -3.91%  '
$cell = $ws.Range('D32')
$cell.NumberFormat = "@"
$cell.Value = '3.03'
$cell.Style = "Normal"
$ws.Range('E32').Value = '  -2.95%  '
$ws.Range('E33').Value = '  -3.97%  '
$cell = $ws.Range('D34')
$cell.NumberFormat = "@"
$cell.Value = '1.46'
$cell.Style = "Normal"
$ws.Range('E34').Value = '  -2.25%  '
$ws.Range('E35').Value = '  -4.10%  '
$ws.Range('D36').Value = '1.087.93'
$ws.Range('E36').Value = '  -2.80%  '
$ws.Range('E37').Value = '  -0.64%  '
$ws.Range('E38').Value = '  -4.82%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$cell = $ws.Range('D40')
$cell.NumberFormat = "@"
$cell.Value = '0.815'
$cell.Style = "Normal"
$ws.Range('E40').Value = '  +6.32%  '
$ws.Range('B41').Value = 'ImmutableX'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$cell = $ws.Range('D41')
$cell.NumberFormat = "@"
$cell.Value = '0.492'
$cell.Style = "Normal"
$ws.Range('E41').Value = '  -4.55%  '
$cell = $ws.Range('D42')
$cell.NumberFormat = "@"
$cell.Value = '0.761'
$cell.Style = "Normal"
$ws.Range('E42').Value = '  -9.94%  '
$cell = $ws.Range('D43')
$cell.NumberFormat = "@"
$cell.Value = '92.98'
$cell.Style = "Normal"
$ws.Range('E43').Value = '  -5.04%  '
$ws.Range('E44').Value = '  -2.01%  '
$ws.Range('D45').Value = '1.693.97'
$ws.Range('E45').Value = '  -3.72%  '
$ws.Range('D46').Value = '0.0₆0112'
$ws.Range('E46').Value = '  -2.95%  '
$cell = $ws.Range('D47')
$cell.NumberFormat = "@"
$cell.Value = '52.45'
$cell.Style = "Normal"
$ws.Range('E47').Value = '  -3.51%  '
$cell = $ws.Range('D48')
$cell.NumberFormat = "@"
$cell.Value = '0.0505'
$cell.Style = "Normal"
$ws.Range('E48').Value = '  -4.98%  '
$ws.Range('E49').Value = '  -2.06%  '
$ws.Range('E50').Value = '  -1.64%  '
$ws.Range('E51').Value = '  -0.34%  '
